$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '57.995.82'
$ws.Range("E2").Value = '  +0.96%  '

# Row 3
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '2.350.91'
$ws.Range("E3").Value = '  +0.98%  '

# Row 4
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '541.26'
$ws.Range("E5").Value = '  -0.19%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '134.78'
$ws.Range("E6").Value = '  -0.09%  '

# Row 7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$ws.Range("E7").Value = '  +0.71%  '

# Row 8
$ws.Range("E8").Value = '  +6.38%  '

# Row 9
$ws.Range("E9").Value = '  +0.66%  '

# Row 10
$ws.Range("E10").Value = '  +2.07%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.357'
$ws.Range("E12").Value = '  +0.92%  '

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '23.82'
$ws.Range("E13").Value = '  +1.19%  '

# Row 14
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '2.766.83'
$ws.Range("E14").Value = '  +0.48%  '

# Row 15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '57.907.39'
$ws.Range("E15").Value = '  +0.64%  '

# Row 16
$ws.Range("E16").Value = '  +0.92%  '

# Row 17
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '2.349.97'
$ws.Range("E17").Value = '  -1.44%  '

# Row 18
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '10.71'
$ws.Range("E18").Value = '  +1.35%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '4.30'
$ws.Range("E19").Value = '  +1.70%  '

# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '330.59'
$ws.Range("E20").Value = '  -2.31%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '6.74'
$ws.Range("E21").Value = '  -1.50%  '

# Row 22
$ws.Range("E22").Value = '  +0.45%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '62.71'
$ws.Range("E23").Value = '  +1.57%  '

# Row 24
$ws.Range("E24").Value = '  -2.41%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$ws.Range("E25").Value = '  -0.05%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '8.36'
$ws.Range("E26").Value = '  -1.79%  '

# Row 27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '1.35'
$ws.Range("E27").Value = '  -5.77%  '

# Row 28
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '170.10'
$ws.Range("E29").Value = '  -0.58%  '

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0737'
$ws.Range("E30").Value = '  -0.04%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '6.13'
$ws.Range("E31").Value = '  -0.78%  '

# Row 32
$ws.Range("E32").Value = '  +1.10%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '18.39'
$ws.Range("E33").Value = '  -0.98%  '

# Row 34
$ws.Range("E34").Value = '  +0.07%  '

# Row 35
$ws.Range("E35").Value = '  +0.92%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '4.21'
$ws.Range("E36").Value = '  +1.28%  '

# Row 37
$ws.Range("E37").Value = '  -1.71%  '

# Row 38
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '1.61'
$ws.Range("E38").Value = '  -0.03%  '

# Row 39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '39.07'
$ws.Range("E39").Value = '  -0.60%  '

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '142.69'
$ws.Range("E40").Value = '  -4.10%  '

# Row 41
$ws.Range("E41").Value = '  -0.02%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '3.65'
$ws.Range("E42").Value = '  +0.53%  '

# Row 43
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '289.19'
$ws.Range("E43").Value = '  +1.28%  '

# Row 44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.0950'
$ws.Range("E44").Value = '  +1.74%  '

# Row 45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.0510'
$ws.Range("E45").Value = '  +0.59%  '

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '19.15'
$ws.Range("E46").Value = '  -0.05%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.566'
$ws.Range("E47").Value = '  +1.25%  '

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.0223'
$ws.Range("E48").Value = '  +1.61%  '

# Row 49
$ws.Range("E49").Value = '  -0.26%  '

# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '11.07'
$ws.Range("E50").Value = '  +0.56%  '

# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '1.55'
$ws.Range("E51").Value = '  -0.65%  '
